# إضافة حدث جديد في Card23 by admin at 2025-12-18 13:05:32
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card23")

# Append the new event as a new row (row 26), by duplicating the previous
# last row (row 25) which already carries the right layout/format
# (card number in A, blank B:K, date/event/correction/servicedBy in L:O).
$ws.Range("A25:O25").Copy($ws.Range("A26:O26"))

# Now that row 25 has been duplicated, fill its previously blank B:K cells
# with "nan" to match the convention used by every other row in the sheet.
$ws.Range("B25:K25").Value = "nan"
